# Technical Design.docx edit script
#
# 1. Removes the spurious <w:proofErr w:type="spellStart"/> / spellEnd
#    markers that wrapped the "Studentnumber" run (the word is not a
#    spelling error, these tags should never have been there).
# 2. Adds a new blank paragraph followed by an "LED driver: ..." paragraph
#    right after the "Technical design: what components specifically"
#    heading paragraph (and before the bookmark end), documenting the
#    interchangeable LED driver ICs.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: drop the proofErr spellStart/spellEnd pair around "Studentnumber"
# ---------------------------------------------------------------------
# Locate the paragraph that starts with "Studentnumber:" and rewrite its
# OOXML without the <w:proofErr/> wrapper elements, keeping the runs
# (and their formatting) identical.
$studentNoPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("Studentnumber:")) {
        $studentNoPara = $candidate
        break
    }
}

$studentNoXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
'<w:body>' +
'<w:p w14:paraId="0E3078B3" w14:textId="606ED6E0" w:rsidR="00A016CA" w:rsidRPr="00615C74" w:rsidRDefault="0011106A" w:rsidP="00D5581E">' +
'<w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Studentnumber</w:t></w:r>' +
'<w:r w:rsidR="00A016CA"><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
'<w:r w:rsidR="00A016CA" w:rsidRPr="00615C74"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>657313</w:t></w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$null = $studentNoPara.Range.InsertXML($studentNoXml)

# ---------------------------------------------------------------------
# Part 2: append a blank paragraph + an "LED driver" paragraph after the
# "Technical design: what components specifically" heading
# ---------------------------------------------------------------------
$headingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith("Technical design: what components")) {
        $headingPara = $candidate
        break
    }
}

# Insert a new (initially empty) paragraph right after the heading; this
# becomes the blank separator paragraph ("<w:p/>" in the target markup).
$null = $headingPara.Range.InsertParagraphAfter()

# Locate the blank paragraph robustly: it is the paragraph immediately
# following the heading paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $headingPara.Range.End) {
        $blankParaIndex = $i
        break
    }
}
$blankPara = $d.Paragraphs($blankParaIndex)

# Insert a second new paragraph after the blank one; this will hold the
# LED driver text.
$null = $blankPara.Range.InsertParagraphAfter()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $blankPara.Range.End) {
        $ledParaIndex = $i
        break
    }
}
$ledPara = $d.Paragraphs($ledParaIndex)

# Make the separator paragraph truly empty (no stray run), matching
# "<w:p/>" in the target OOXML.
$blankEmptyXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body><w:p/></w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'
$null = $blankPara.Range.InsertXML($blankEmptyXml)

# Fill the LED driver paragraph with the required runs.
$ledXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:r><w:t xml:space="preserve">LED driver: </w:t></w:r>' +
'<w:r><w:t>TPS61169</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> or </w:t></w:r>' +
'<w:r><w:t>TPS92360</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> (interchangeable)</w:t></w:r>' +
'</w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'
$null = $ledPara.Range.InsertXML($ledXml)

Write-Output "Edit complete."
